$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K5").Value = -0.3802
$ws.Range("L5").Value = -0.1156
$ws.Range("M5").Value = -0.0439
$ws.Range("N5").Value = -0.2949
$ws.Range("O5").Value = -0.4257
$ws.Range("P5").Value = -0.0613
$ws.Range("Q5").Value = -0.0515
$ws.Range("R5").Value = -0.3876
$ws.Range("S5").Value = -32.9216
$ws.Range("K7").Value = 0.0179
$ws.Range("L7").Value = 0.0029
$ws.Range("M7").Value = -0.0022
$ws.Range("O7").Value = -0.0113
$ws.Range("P7").Value = 0.007
$ws.Range("Q7").Value = -0.0111
$ws.Range("R7").Value = -0.0244
$ws.Range("S7").Value = -5.4036
$ws.Range("K8").Value = 0.3134
$ws.Range("L8").Value = 0.2152
$ws.Range("M8").Value = 0.2641
$ws.Range("N8").Value = 0.5368
$ws.Range("O8").Value = 0.5283
$ws.Range("P8").Value = 0.5569
$ws.Range("Q8").Value = 0.5092
$ws.Range("R8").Value = 0.2708
$ws.Range("S8").Value = 7.6505
$ws.Range("K12").Value = 0.225
$ws.Range("L12").Value = 0.2688
$ws.Range("M12").Value = -0.0137
$ws.Range("N12").Value = -0.1003
$ws.Range("O12").Value = -0.0852
$ws.Range("P12").Value = -0.0702
$ws.Range("Q12").Value = -0.0351
$ws.Range("R12").Value = -0.0729
$ws.Range("S12").Value = -6.6924
$ws.Range("K16").Value = -1.2456
$ws.Range("L16").Value = -0.1857
$ws.Range("M16").Value = -0.1632
$ws.Range("N16").Value = -0.3603
$ws.Range("O16").Value = 0.0317
$ws.Range("P16").Value = 0.0696
$ws.Range("Q16").Value = 0.1879
$ws.Range("R16").Value = -0.0504
$ws.Range("S16").Value = -72.2455
$ws.Range("K18").Value = -0.1144
$ws.Range("L18").Value = -0.1186
$ws.Range("M18").Value = -0.0684
$ws.Range("N18").Value = -0.2598
$ws.Range("O18").Value = -0.3725
$ws.Range("P18").Value = -0.0132
$ws.Range("Q18").Value = 0.0133
$ws.Range("R18").Value = -0.3159
$ws.Range("S18").Value = -6.8434
$ws.Range("K22").Value = 0.0681
$ws.Range("L22").Value = 0.123
$ws.Range("M22").Value = -0.0103
$ws.Range("N22").Value = 0.2002
$ws.Range("O22").Value = 0.2917
$ws.Range("P22").Value = -0.0699
$ws.Range("Q22").Value = -0.0998
$ws.Range("R22").Value = 0.2271
$ws.Range("S22").Value = -37.433
$ws.Range("K24").Value = 0.0149
$ws.Range("L24").Value = 0.1232
$ws.Range("M24").Value = 0.0766
$ws.Range("O24").Value = -0.0182
$ws.Range("P24").Value = -0.0672
$ws.Range("Q24").Value = -0.0456
$ws.Range("R24").Value = -0.036
$ws.Range("S24").Value = -0.7106
$ws.Range("K29").Value = -0.0131
$ws.Range("L29").Value = 0.006
$ws.Range("M29").Value = 0.0207
$ws.Range("N29").Value = 0.0219
$ws.Range("O29").Value = 0.0242
$ws.Range("P29").Value = 0.0184
$ws.Range("Q29").Value = 0.0086
$ws.Range("R29").Value = 0.0021
$ws.Range("S29").Value = -0.0885
$ws.Range("K33").Value = -0.33
$ws.Range("L33").Value = -0.0631
$ws.Range("M33").Value = 0.0005
$ws.Range("N33").Value = 0.0012
$ws.Range("O33").Value = 0.0022
$ws.Range("P33").Value = 0.0005
$ws.Range("Q33").Value = -0.0011
$ws.Range("R33").Value = 0.0029
$ws.Range("S33").Value = -0.0709
$ws.Range("K35").Value = 0.0163
$ws.Range("L35").Value = 0.0449
$ws.Range("M35").Value = 0.0392
$ws.Range("N35").Value = -0.0216
$ws.Range("O35").Value = 0.0127
$ws.Range("P35").Value = -0.0197
$ws.Range("Q35").Value = -0.042
$ws.Range("R35").Value = -0.0654
$ws.Range("S35").Value = -0.0003
$ws.Range("K36").Value = 0.3056
$ws.Range("L36").Value = 0.1527
$ws.Range("M36").Value = 0.139
$ws.Range("N36").Value = 0.1345
$ws.Range("O36").Value = 0.1312
$ws.Range("P36").Value = 0.1275
$ws.Range("Q36").Value = 0.127
$ws.Range("R36").Value = -0.0085
$ws.Range("S36").Value = -0.3384
$ws.Range("K40").Value = 0.0992
$ws.Range("L40").Value = 0.0969
$ws.Range("M40").Value = 0.0967
$ws.Range("N40").Value = -0.0046
$ws.Range("O40").Value = -0.0044
$ws.Range("P40").Value = -0.0042
$ws.Range("Q40").Value = -0.0042
$ws.Range("R40").Value = -0.0037
$ws.Range("S40").Value = -0.0608
$ws.Range("K44").Value = 0.3152
$ws.Range("L44").Value = 0.2717
$ws.Range("M44").Value = 0.3192
$ws.Range("N44").Value = 0.182
$ws.Range("O44").Value = 0.1856
$ws.Range("P44").Value = 0.1827
$ws.Range("Q44").Value = 0.1593
$ws.Range("R44").Value = -0.0568
$ws.Range("S44").Value = -0.1023
$ws.Range("K46").Value = -0.083
$ws.Range("L46").Value = 0.0003
$ws.Range("M46").Value = 0
$ws.Range("N46").Value = 0.0014
$ws.Range("O46").Value = 0.0023
$ws.Range("P46").Value = -0.0001
$ws.Range("Q46").Value = -0.0003
$ws.Range("R46").Value = 0.0019
$ws.Range("S46").Value = -0.1259
$ws.Range("K50").Value = 0.0834
$ws.Range("L50").Value = -0.0003
$ws.Range("M50").Value = 0.0007
$ws.Range("N50").Value = -0.0009
$ws.Range("O50").Value = -0.0017
$ws.Range("P50").Value = 0.0008
$ws.Range("Q50").Value = 0.001
$ws.Range("R50").Value = -0.0012
$ws.Range("S50").Value = 0.4949
$ws.Range("K52").Value = 0.0093
$ws.Range("L52").Value = -0.0267
$ws.Range("M52").Value = -0.0272
$ws.Range("O52").Value = -0.0369
$ws.Range("P52").Value = -0.0032
$ws.Range("Q52").Value = 0.0052
$ws.Range("R52").Value = 0.0132
$ws.Range("S52").Value = 0.0324
$ws.Range("K57").Value = -0.0069
$ws.Range("L57").Value = -0.001
$ws.Range("M57").Value = 0.0047
$ws.Range("N57").Value = 0.0096
$ws.Range("O57").Value = 0.0167
$ws.Range("P57").Value = 0.0159
$ws.Range("Q57").Value = 0.0102
$ws.Range("R57").Value = 0.0061
$ws.Range("S57").Value = -0.009
